# CropLifeCycle.xlsx update
# Commit: "Updated Crop Life Cycle and Model Plan"
#   Changed sowing month of Paddy from 1st June to 1st July and
#   Similarly Harvest month from 1st September to 1st October.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 of the "Table1" structured table is the Paddy/Rice crop entry:
#   A8=Crop, B8=Season, C8=Sowing Start, D8=Sowing End,
#   E8=Harvest Start, F8=Harvest End, G8=Average Months
# Update the sowing start date (C8) and harvesting start date (E8),
# and the recalculated average months figure (G8).
$ws.Range("C8").Value = "1st July"
$ws.Range("E8").Value = "1st October"
$ws.Range("G8").Value = 5

# Reflect the active cell/selection left behind after the edit.
[void]$ws.Range("G13").Select()
